# Splits the original "Abc xyz" run into spell-checked sub-runs (as Word's
# proofing pass would after an edit), then appends a blank paragraph and a
# new paragraph containing "Master branch new changes ".

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Rewrite paragraph 1 ("Abc xyz") into proofErr-wrapped runs:
#    Abc | " " | xyz, each word flagged with spellStart/spellEnd.
$p1 = $d.Paragraphs(1)
$firstXml = '<w:p ' + $wNs + '>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>Abc</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
              '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:proofErr w:type="spellStart"/>' +
              '<w:r><w:t>xyz</w:t></w:r>' +
              '<w:proofErr w:type="spellEnd"/>' +
            '</w:p>'
[void]$p1.Range.InsertXML($firstXml)

# 2) Append a blank paragraph and the "Master branch new changes " paragraph
#    right after it (before the end-of-story mark).
$endPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newXml = '<w:p ' + $wNs + '/>' +
          '<w:p ' + $wNs + '>' +
            '<w:r><w:t xml:space="preserve">Master branch new changes </w:t></w:r>' +
          '</w:p>'
[void]$endPoint.InsertXML($newXml)
